# Apply updated Mann-Kendall duration-calc outputs and derived site_metrics trend flags.
$wb = $excel.ActiveWorkbook

# site_metrics
$ws1 = $wb.Worksheets.Item("site_metrics")
$ws1.Range("AK13").Value = $true
$ws1.Range("AK14").Value = $true
$ws1.Range("O21").Value = 0.1559653515267294
$ws1.Range("O26").Value = 0.0006684848531900506
$ws1.Range("N32").Value = 5.1322175840033
$ws1.Range("O32").Value = 0.01157511800627542
$ws1.Range("O33").Value = 0.03040019100535623
$ws1.Range("AK36").Value = $true
$ws1.Range("AK37").Value = $true
$ws1.Range("AK40").Value = $true
$ws1.Range("AK41").Value = $true
$ws1.Range("O43").Value = 0.005877892448927921
$ws1.Range("AK43").Value = $true
$ws1.Range("AK44").Value = $true
$ws1.Range("O49").Value = 0.01430634299687951
$ws1.Range("O50").Value = 0.01783381014182999
$ws1.Range("AK52").Value = $true
$ws1.Range("AK60").Value = $true
$ws1.Range("O65").Value = 0.00516518207657336
$ws1.Range("O66").Value = 0.005663890535828419
$ws1.Range("AK69").Value = $true
$ws1.Range("Q70").Value = 2.825
$ws1.Range("N71").Value = 3.562547251223723
$ws1.Range("O71").Value = 0.003729460387189254
$ws1.Range("AK73").Value = $true
$ws1.Range("N75").Value = 5.498488031576267
$ws1.Range("O75").Value = 0.001274927373175842
$ws1.Range("O79").Value = 0.02058574564297492

# mk_duration
$ws3 = $wb.Worksheets.Item("mk_duration")
$ws3.Range("K5").Value = "no trend"
$ws3.Range("L5").Value = $false
$ws3.Range("M5").Value = 0.05016562050620044
$ws3.Range("N5").Value = -1.958549057513455
$ws3.Range("O5").Value = -0.2586206896551724
$ws3.Range("P5").Value = -105
$ws3.Range("Q5").Value = 2819.666666666667
$ws3.Range("R5").Value = -0.3100961538461539
$ws3.Range("S5").Value = 8.341346153846153
$ws3.Range("M26").Value = 0.3363974735521627
$ws3.Range("N26").Value = -0.9613077077685092
$ws3.Range("O26").Value = -0.1193181818181818
$ws3.Range("P26").Value = -63
$ws3.Range("Q26").Value = 4159.666666666667
$ws3.Range("R26").Value = -0.008319887784173498
$ws3.Range("S26").Value = 2.258118204546776
$ws3.Range("M31").Value = 0.7915094652925039
$ws3.Range("N31").Value = -0.2643509998273594
$ws3.Range("O31").Value = -0.02570921985815603
$ws3.Range("P31").Value = -29
$ws3.Range("Q31").Value = 11219
$ws3.Range("S31").Value = 1
$ws3.Range("K32").Value = "no trend"
$ws3.Range("L32").Value = $false
$ws3.Range("M32").Value = 0.3373746480667339
$ws3.Range("N32").Value = -0.9593655015712707
$ws3.Range("O32").Value = -0.1280788177339902
$ws3.Range("P32").Value = -52
$ws3.Range("Q32").Value = 2826
$ws3.Range("R32").Value = -0.05941876750700281
$ws3.Range("S32").Value = 5.94297385620915
$ws3.Range("K39").Value = "no trend"
$ws3.Range("L39").Value = $false
$ws3.Range("M39").Value = 0.3519057100306646
$ws3.Range("N39").Value = 0.9308991968322524
$ws3.Range("O39").Value = 0.13
$ws3.Range("P39").Value = 39
$ws3.Range("Q39").Value = 1666.333333333333
$ws3.Range("S39").Value = 4
$ws3.Range("M61").Value = 0.8095827631119468
$ws3.Range("N61").Value = 0.2409643292909423
$ws3.Range("O61").Value = 0.026578073089701
$ws3.Range("P61").Value = 24
$ws3.Range("Q61").Value = 9110.666666666666
$ws3.Range("R61").Value = 0.007017543859649121
$ws3.Range("S61").Value = 5.519298245614036
$ws3.Range("M65").Value = 0.8277368747621565
$ws3.Range("N65").Value = -0.2176050175140125
$ws3.Range("O65").Value = -0.03174603174603174
$ws3.Range("P65").Value = -12
$ws3.Range("Q65").Value = 2555.333333333333
$ws3.Range("R65").Value = -0.01481481481481482
$ws3.Range("S65").Value = 6
$ws3.Range("M66").Value = 0.1463622442996033
$ws3.Range("N66").Value = -1.452501385187843
$ws3.Range("O66").Value = -0.1681681681681682
$ws3.Range("P66").Value = -112
$ws3.Range("Q66").Value = 5840
$ws3.Range("R66").Value = -0.1060139573070608
$ws3.Range("S66").Value = 8.283251231527094
$ws3.Range("D70").Value = 0.0002126871922845108
$ws3.Range("E70").Value = -3.703449909656969
$ws3.Range("F70").Value = -0.4143049932523617
$ws3.Range("G70").Value = -307
$ws3.Range("H70").Value = 6827
$ws3.Range("I70").Value = -0.5217391304347826
$ws3.Range("J70").Value = 22.16304347826087
$ws3.Range("M70").Value = 0.02065253076358409
$ws3.Range("N70").Value = -2.314277033323547
$ws3.Range("O70").Value = -0.2253061224489796
$ws3.Range("P70").Value = -276
$ws3.Range("Q70").Value = 14120
$ws3.Range("R70").Value = -0.25
$ws3.Range("S70").Value = 13.875
$ws3.Range("M71").Value = 0.9170221576248359
$ws3.Range("N71").Value = 0.104185478635438
$ws3.Range("O71").Value = 0.01159420289855072
$ws3.Range("P71").Value = 12
$ws3.Range("Q71").Value = 11147.33333333333
$ws3.Range("R71").Value = 0.001190476190476186
$ws3.Range("S71").Value = 2.540178571428571
$ws3.Range("M75").Value = 0.4056376702457534
$ws3.Range("N75").Value = -0.8315948495541449
$ws3.Range("O75").Value = -0.09841269841269841
$ws3.Range("P75").Value = -62
$ws3.Range("Q75").Value = 5380.666666666667
$ws3.Range("R75").Value = -0.04517770876466531
$ws3.Range("S75").Value = 4.665609903381643
$ws3.Range("M76").Value = 0.7191153342667906
$ws3.Range("N76").Value = -0.3596413803295871
$ws3.Range("O76").Value = -0.04024390243902439
$ws3.Range("P76").Value = -33
$ws3.Range("Q76").Value = 7917
$ws3.Range("R76").Value = -0.0109775641025641
$ws3.Range("S76").Value = 2.552884615384615
$ws3.Range("M78").Value = 0.7981372539669001
$ws3.Range("N78").Value = -0.2557585861466805
$ws3.Range("O78").Value = -0.02705314009661836
$ws3.Range("P78").Value = -28
$ws3.Range("Q78").Value = 11144.66666666667
$ws3.Range("R78").Value = -0.007352941176470588
$ws3.Range("S78").Value = 6.415441176470588
$ws3.Range("M80").Value = 0.8304218583266509
$ws3.Range("N80").Value = 0.2141605507090119
$ws3.Range("O80").Value = 0.02988505747126437
$ws3.Range("P80").Value = 13
$ws3.Range("Q80").Value = 3139.666666666667
$ws3.Range("R80").Value = 0.0166666666666667
$ws3.Range("S80").Value = 5.283333333333333
$ws3.Range("M82").Value = 0.5236941825672736
$ws3.Range("N82").Value = 0.6376613001696565
$ws3.Range("O82").Value = 0.07843137254901961
$ws3.Range("P82").Value = 44
$ws3.Range("Q82").Value = 4547.333333333333
$ws3.Range("R82").Value = 0.03529411764705883
$ws3.Range("S82").Value = 3.7390756302521

# mk_intra_annual
$ws4 = $wb.Worksheets.Item("mk_intra_annual")
$ws4.Range("K5").Value = "no trend"
$ws4.Range("L5").Value = $false
$ws4.Range("M5").Value = 0.3808184731875461
$ws4.Range("N5").Value = -0.8763892264428205
$ws4.Range("O5").Value = -0.1157635467980296
$ws4.Range("P5").Value = -47
$ws4.Range("Q5").Value = 2755
$ws4.Range("S5").Value = 3
$ws4.Range("K26").Value = "no trend"
$ws4.Range("L26").Value = $false
$ws4.Range("M26").Value = 0.4183732951800732
$ws4.Range("N26").Value = -0.8092466349951097
$ws4.Range("O26").Value = -0.1003787878787879
$ws4.Range("P26").Value = -53
$ws4.Range("Q26").Value = 4129
$ws4.Range("R26").Value = -0.07692307692307693
$ws4.Range("S26").Value = 17.23076923076923
$ws4.Range("M31").Value = 0.1711961397741817
$ws4.Range("N31").Value = -1.368370479444599
$ws4.Range("O31").Value = -0.1285460992907801
$ws4.Range("P31").Value = -145
$ws4.Range("Q31").Value = 11074.33333333333
$ws4.Range("S31").Value = 1
$ws4.Range("K32").Value = "no trend"
$ws4.Range("L32").Value = $false
$ws4.Range("M32").Value = 0.6072965855556409
$ws4.Range("N32").Value = 0.5139362166481446
$ws4.Range("O32").Value = 0.06896551724137931
$ws4.Range("P32").Value = 28
$ws4.Range("Q32").Value = 2760
$ws4.Range("R32").Value = 0
$ws4.Range("S32").Value = 6
$ws4.Range("K39").Value = "no trend"
$ws4.Range("L39").Value = $false
$ws4.Range("M39").Value = 0.3306659971936292
$ws4.Range("N39").Value = 0.9727732764453586
$ws4.Range("O39").Value = 0.1333333333333333
$ws4.Range("P39").Value = 40
$ws4.Range("Q39").Value = 1607.333333333333
$ws4.Range("S39").Value = 1
$ws4.Range("M61").Value = 0.4166463826930014
$ws4.Range("N61").Value = -0.8122531581214926
$ws4.Range("O61").Value = -0.08637873754152824
$ws4.Range("P61").Value = -78
$ws4.Range("Q61").Value = 8986.666666666666
$ws4.Range("S61").Value = 5
$ws4.Range("M65").Value = 0.9522762374920219
$ws4.Range("N65").Value = -0.05984857517200719
$ws4.Range("O65").Value = -0.01058201058201058
$ws4.Range("P65").Value = -4
$ws4.Range("Q65").Value = 2512.666666666667
$ws4.Range("S65").Value = 4.5
$ws4.Range("M66").Value = 0.9894548608294731
$ws4.Range("N66").Value = -0.01321675678200343
$ws4.Range("O66").Value = -0.003003003003003003
$ws4.Range("P66").Value = -2
$ws4.Range("Q66").Value = 5724.666666666667
$ws4.Range("R66").Value = 0
$ws4.Range("S66").Value = 5
$ws4.Range("M71").Value = 0.005350557452794114
$ws4.Range("N71").Value = 2.785134192734321
$ws4.Range("O71").Value = 0.2840579710144928
$ws4.Range("P71").Value = 294
$ws4.Range("Q71").Value = 11067.33333333333
$ws4.Range("R71").Value = 0.1428571428571428
$ws4.Range("S71").Value = 6.785714285714286
$ws4.Range("M75").Value = 0.6021588726858282
$ws4.Range("N75").Value = 0.5212984613601388
$ws4.Range("O75").Value = 0.06190476190476191
$ws4.Range("P75").Value = 39
$ws4.Range("Q75").Value = 5313.666666666667
$ws4.Range("R75").Value = 0
$ws4.Range("S75").Value = 6.5
$ws4.Range("K76").Value = "no trend"
$ws4.Range("L76").Value = $false
$ws4.Range("M76").Value = 0.1087367540347022
$ws4.Range("N76").Value = 1.603897018969444
$ws4.Range("O76").Value = 0.174390243902439
$ws4.Range("P76").Value = 143
$ws4.Range("Q76").Value = 7838.333333333333
$ws4.Range("R76").Value = 0.1096096096096096
$ws4.Range("S76").Value = 7.807807807807808
$ws4.Range("M78").Value = 0.6448044792368861
$ws4.Range("N78").Value = -0.4609918124993016
$ws4.Range("O78").Value = -0.04734299516908213
$ws4.Range("P78").Value = -49
$ws4.Range("Q78").Value = 10841.66666666667
$ws4.Range("M80").Value = 0.5993888518643082
$ws4.Range("N80").Value = -0.5252795781016084
$ws4.Range("O80").Value = -0.06896551724137931
$ws4.Range("P80").Value = -30
$ws4.Range("Q80").Value = 3048
$ws4.Range("S80").Value = 6
$ws4.Range("K82").Value = "no trend"
$ws4.Range("L82").Value = $false
$ws4.Range("M82").Value = 0.2881499292551952
$ws4.Range("N82").Value = -1.062188919124312
$ws4.Range("O82").Value = -0.1283422459893048
$ws4.Range("P82").Value = -72
$ws4.Range("Q82").Value = 4468
$ws4.Range("R82").Value = -0.05263157894736842
$ws4.Range("S82").Value = 6.868421052631579
